# edit.ps1
# Applies the "Updated symbol list" GitHub Actions refresh (Tue Dec 20
# 10:48:45 UTC 2022): new Price quotes for most coins, a couple of rows
# whose coin/link/volume-label data got reshuffled (KickToken / BKEXToken /
# CEJI rotated down one row), and one "Worst in 24h" label update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store numeric-looking values as TEXT in this
# workbook. A leading apostrophe keeps Range.Value from auto-converting the
# new quotes to numbers, matching the original text-typed cells.
function Set-TextValue($range, [string]$value) {
    $ws.Range($range).Value = "'" + $value
}

# Row 2
Set-TextValue "D2" "248.22"

# Row 3
Set-TextValue "D3" "21.96"

# Row 4
Set-TextValue "D4" "5.366"

# Row 5
Set-TextValue "D5" "0.05632"

# Row 6
Set-TextValue "D6" "3.434"

# Row 7
Set-TextValue "D7" "6.350"

# Row 8
Set-TextValue "D8" "0.8161"

# Row 9
Set-TextValue "D9" "0.9268"

# Row 10
Set-TextValue "D10" "0.1430"

# Row 11
Set-TextValue "D11" "0.07450"

# Row 12
Set-TextValue "D12" "0.03255"

# Row 13
Set-TextValue "D13" "0.03094"

# Row 14
Set-TextValue "D14" "0.09329"

# Row 15
Set-TextValue "D15" "3.558"

# Row 16
Set-TextValue "D16" "0.001596"

# Row 17
Set-TextValue "D17" "0.04720"

# Row 18
Set-TextValue "D18" "0.0005785"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19
Set-TextValue "D19" "0.006396"

# Row 20
Set-TextValue "D20" "0.005054"

# Row 23
Set-TextValue "D23" "3.745"

# Row 26
Set-TextValue "D26" "0.1320"

# Row 40
Set-TextValue "D40" "0.03942"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006874"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003402"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 45
Set-TextValue "D45" "0.00005577"

# Row 48
Set-TextValue "D48" "0.7806"

# Row 49
Set-TextValue "D49" "0.1799"

# Row 51
Set-TextValue "D51" "0.01011"

